# Update naive QoQ matched-to-ifoCAST error tables with full-series evaluation values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1666076251416388
$ws.Range("C2").Value = 1.464488876860186
$ws.Range("D2").Value = 9.385379659755237
$ws.Range("E2").Value = 3.063556700920556
$ws.Range("F2").Value = 3.131009834800738
$ws.Range("G2").Value = 22

$ws.Range("B3").Value = 0.0708499967161877
$ws.Range("C3").Value = 1.419337315890072
$ws.Range("D3").Value = 9.02840922660916
$ws.Range("E3").Value = 3.004731140486476
$ws.Range("F3").Value = 3.078077156245957
$ws.Range("G3").Value = 21

$ws.Range("B4").Value = -0.4288746198142434
$ws.Range("C4").Value = 0.8403831407438848
$ws.Range("D4").Value = 3.678417349840328
$ws.Range("E4").Value = 1.917920058250689
$ws.Range("F4").Value = 1.917916608230349
$ws.Range("G4").Value = 20

$ws.Range("B5").Value = -0.01747259896522734
$ws.Range("C5").Value = 0.5856938463863763
$ws.Range("D5").Value = 0.8820486078156717
$ws.Range("E5").Value = 0.9391744288552961
$ws.Range("F5").Value = 0.9647429999147259
$ws.Range("G5").Value = 19

$ws.Range("B6").Value = -0.01422969815280487
$ws.Range("C6").Value = 0.6208678124689757
$ws.Range("D6").Value = 0.9450450342577459
$ws.Range("E6").Value = 0.9721342676080017
$ws.Range("F6").Value = 1.000210739531721
$ws.Range("G6").Value = 18

$ws.Range("B7").Value = -0.1197329514763359
$ws.Range("C7").Value = 0.4398200131717107
$ws.Range("D7").Value = 0.4468041506909826
$ws.Range("E7").Value = 0.6684341034769117
$ws.Range("F7").Value = 0.6778623988027417
$ws.Range("G7").Value = 17

$ws.Range("B8").Value = -0.04262487491856422
$ws.Range("C8").Value = 0.4573970767948405
$ws.Range("D8").Value = 0.4014438679027904
$ws.Range("E8").Value = 0.6335959816024643
$ws.Range("F8").Value = 0.6528926306346485
$ws.Range("G8").Value = 16

$ws.Range("B9").Value = 0.04839913305653263
$ws.Range("C9").Value = 0.3941245839882653
$ws.Range("D9").Value = 0.3208696055700829
$ws.Range("E9").Value = 0.566453533460674
$ws.Range("F9").Value = 0.5841909510683425
$ws.Range("G9").Value = 15

$ws.Range("B10").Value = -0.01748806567405372
$ws.Range("C10").Value = 0.3260429352179946
$ws.Range("D10").Value = 0.2494121353768383
$ws.Range("E10").Value = 0.4994117893851108
$ws.Range("F10").Value = 0.5179462580602099
$ws.Range("G10").Value = 14

$ws.Range("B11").Value = -0.02168945670833724
$ws.Range("C11").Value = 0.351174191232308
$ws.Range("D11").Value = 0.2265631253025311
$ws.Range("E11").Value = 0.4759864759659995
$ws.Range("F11").Value = 0.4949078202060265
$ws.Range("G11").Value = 13

